# Resume update for websphere move
$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1, wdReplaceAll = 2

# 1. "Programming Languages" list: move "Go" from the end to the front.
$d.Content.Find.Execute(
    "Python, Perl, C#, C++, Go", $false, $false, $false, $false, $false,
    $true, 1, $false, "Go, Python, Perl, C#, C++", 2)

# 2. "Technologies/Frameworks" list: add Swagger, Cobra right after the colon.
$d.Content.Find.Execute(
    ": VMware NSX,", $false, $false, $false, $false, $false,
    $true, 1, $false, ": Swagger, Cobra, VMware NSX,", 2)

# 3. First EXPERIENCE heading "IBM: WebSphere " -> "IBM " (only the first
#    occurrence - the second "IBM: WebSphere " entry, for the earlier
#    internship, is left alone).
$d.Content.Find.Execute(
    "IBM: WebSphere ", $false, $false, $false, $false, $false,
    $true, 1, $false, "IBM ", 1)

# 4. Job title line loses the "| WebSphere as a Service, Networking" suffix
#    (that text moves down into the bullet below).
$d.Content.Find.Execute(
    "Software Developer | WebSphere as a Service, Networking",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "Software Developer ", 1)

# 5. First bullet now starts by naming the team/project, and "Manage" is
#    lower-cased since it is no longer the first word of the sentence.
$d.Content.Find.Execute(
    "Manage and automate VMware network architecture for WASaas environments using Python and Ansible.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "WebSphere as a Service, Networking: manage and automate VMware network architecture for WASaaS environments using Python and Ansible.",
    1)

# 6. Add a new bullet describing the Kabanero CLI work, right after the bullet
#    we just edited (paragraph 17), reusing that paragraph's bullet/list
#    formatting.
$bulletPara = $d.Paragraphs(17)
$bulletPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(18)
$newRange = $newPara.Range
$newRange.InsertAfter("Build the ")
$newRange.InsertAfter("Kabanero")
$newRange.InsertAfter(" CLI in Go to handle stack management for governed application development.")

# Move the "_GoBack" bookmark (Word's "last edit location" marker) from its
# old spot near "(2018)." to the middle of the sentence we just typed, which
# is where the document's last edit now is.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$goBackAnchor = $d.Content
$goBackAnchor.Find.Execute("governe")
$goBackPos = $goBackAnchor.End
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

Write-Output "done"
